$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.513.38'
$ws.Range("E2").Value = '  +2.53%  '

$ws.Range("D3").Value = '2.312.64'
$ws.Range("E3").Value = '  +1.68%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''311.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '

$ws.Range("D6").Value = '''102.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.87%  '

$ws.Range("E7").Value = '  +1.55%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '''0.532'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.13%  '

$ws.Range("D10").Value = '''35.76'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("D11").Value = '''0.0818'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.22%  '

$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("E13").Value = '  +1.02%  '

$ws.Range("D14").Value = '2.669.77'

$ws.Range("D15").Value = '''14.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").Value = '2.310.69'
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("E17").Value = '  +1.66%  '

$ws.Range("D18").Value = '43.416.61'
$ws.Range("E18").Value = '  +2.66%  '

$ws.Range("D19").Value = '''12.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.91%  '

$ws.Range("D20").Value = '0.0₃0927'
$ws.Range("E20").Value = '  +2.08%  '

$ws.Range("D21").Value = '''6.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.29%  '

$ws.Range("D22").Value = '''68.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = '''241.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.13%  '

$ws.Range("D24").Value = '''2.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.24%  '

$ws.Range("D25").Value = '''2.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.46%  '

$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("E27").Value = '  -1.60%  '

$ws.Range("D28").Value = '''24.83'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.75%  '

$ws.Range("D29").Value = '''2.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.24%  '

$ws.Range("D30").Value = '''36.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.68%  '

$ws.Range("D31").Value = '''9.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.71%  '

$ws.Range("D32").Value = '''167.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.70%  '

$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.0746'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.44%  '

$ws.Range("E36").Value = '  -2.61%  '

$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '''2.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.68%  '

$ws.Range("D38").Value = '''17.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.87%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.106'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.28%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '''1.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.23%  '

$ws.Range("E41").Value = '  +1.38%  '

$ws.Range("D42").Value = '''4.32'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.51%  '

$ws.Range("E43").Value = '  -0.75%  '

$ws.Range("D44").Value = '''19.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.44%  '

$ws.Range("E45").Value = '  +2.05%  '

$ws.Range("D46").Value = '1.970.10'
$ws.Range("E46").Value = '  +1.10%  '

$ws.Range("E47").Value = '  +2.30%  '

$ws.Range("D48").Value = '''9.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.05%  '

$ws.Range("D49").Value = '''55.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.09%  '

$ws.Range("E50").Value = '  +5.75%  '

$ws.Range("E51").Value = '  +6.00%  '
